$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H63").Value = 127499.75
$ws.Range("J63").Value = 127499.75
$ws.Range("L63").Value = 127499.75
$ws.Range("N63").Value = -128747.75
$ws.Range("H66").Value = 127499.75
$ws.Range("J66").Value = 127499.75
$ws.Range("L66").Value = 382499.25
$ws.Range("N66").Value = -388739.25
$ws.Range("H80").Value = 1044.6207
$ws.Range("I80").Value = 665.38464
$ws.Range("K80").Value = 1996.15392
$ws.Range("M80").Value = -998.15392
$ws.Range("H83").Value = 1044.6207
$ws.Range("I83").Value = 665.38464
$ws.Range("K83").Value = 5988.46176
$ws.Range("M83").Value = -996.4617600000001
$ws.Range("H115").Value = 637.4
$ws.Range("I115").Value = 637.4
$ws.Range("K115").Value = 1912.2
$ws.Range("M115").Value = -345.1999999999998
$ws.Range("H132").Value = 22746.586
$ws.Range("I132").Value = 3907.25
$ws.Range("K132").Value = 11721.75
$ws.Range("M132").Value = -9191.75
$ws.Range("H138").Value = 3594.2395
$ws.Range("I138").Value = 1851.7273
$ws.Range("J138").Value = 4376.592
$ws.Range("K138").Value = 5555.1819
$ws.Range("L138").Value = 13129.776
$ws.Range("M138").Value = -415.1818999999996
$ws.Range("N138").Value = -23409.776

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 2452.75
$ws.Range("I2").Value = 2672.5625
$ws.Range("J2").Value = 1573.5
$ws.Range("K2").Value = 2672.5625
$ws.Range("L2").Value = 1573.5
$ws.Range("M2").Value = -2559.5625
$ws.Range("N2").Value = -1799.5
$ws.Range("H32").Value = 12663854
$ws.Range("I32").Value = 13163506
$ws.Range("K32").Value = 13163506
$ws.Range("M32").Value = -13163219
$ws.Range("H61").Value = 2283.3333
$ws.Range("I61").Value = 2031.2307
$ws.Range("J61").Value = 2938.8
$ws.Range("K61").Value = 2031.2307
$ws.Range("L61").Value = 2938.8
$ws.Range("M61").Value = -1819.2307
$ws.Range("N61").Value = -3362.8
$ws.Range("H116").Value = 2452.75
$ws.Range("I116").Value = 2672.5625
$ws.Range("J116").Value = 1573.5
$ws.Range("K116").Value = 2672.5625
$ws.Range("L116").Value = 1573.5
$ws.Range("M116").Value = -378.5625
$ws.Range("N116").Value = -6161.5
$ws.Range("H132").Value = 3027.1428
$ws.Range("I132").Value = 3027.1428
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9081.428400000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6551.428400000001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2283.3333
$ws.Range("I136").Value = 2031.2307
$ws.Range("J136").Value = 2938.8
$ws.Range("K136").Value = 6093.6921
$ws.Range("L136").Value = 8816.400000000001
$ws.Range("M136").Value = -3543.6921
$ws.Range("N136").Value = -13916.4

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 2452.75
$ws.Range("I3").Value = 2672.5625
$ws.Range("J3").Value = 1573.5
$ws.Range("K3").Value = 2672.5625
$ws.Range("L3").Value = 1573.5
$ws.Range("M3").Value = -2558.5625
$ws.Range("N3").Value = -1801.5
$ws.Range("H86").Value = 2974.9211
$ws.Range("I86").Value = 2217.6538
$ws.Range("J86").Value = 4615.6665
$ws.Range("K86").Value = 2217.6538
$ws.Range("L86").Value = 4615.6665
$ws.Range("M86").Value = -1094.6538
$ws.Range("N86").Value = -6861.6665
$ws.Range("H89").Value = 2974.9211
$ws.Range("I89").Value = 2217.6538
$ws.Range("J89").Value = 4615.6665
$ws.Range("K89").Value = 11088.269
$ws.Range("L89").Value = 23078.3325
$ws.Range("M89").Value = -5472.269
$ws.Range("N89").Value = -34310.3325
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 10410.9375
$ws.Range("I107").Value = 8198.076999999999
$ws.Range("K107").Value = 8198.076999999999
$ws.Range("M107").Value = -6278.076999999999
$ws.Range("H134").Value = 1729.55
$ws.Range("I134").Value = 1348.0571
$ws.Range("K134").Value = 4044.1713
$ws.Range("M134").Value = -1509.1713
$ws.Range("H140").Value = 73984.07000000001

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H22").Value = 602.9
$ws.Range("I22").Value = 436.6
$ws.Range("J22").Value = 769.2
$ws.Range("K22").Value = 436.6
$ws.Range("L22").Value = 769.2
$ws.Range("M22").Value = -86.60000000000002
$ws.Range("N22").Value = -1469.2
$ws.Range("H31").Value = 2034.025
$ws.Range("I31").Value = 1966.6897
$ws.Range("K31").Value = 1966.6897
$ws.Range("M31").Value = -1671.6897
$ws.Range("H34").Value = 2034.025
$ws.Range("I34").Value = 1966.6897
$ws.Range("K34").Value = 1966.6897
$ws.Range("M34").Value = -1764.6897
$ws.Range("H122").Value = 791478.25
$ws.Range("I122").Value = 2554325.8
$ws.Range("J122").Value = 7990.4443
$ws.Range("K122").Value = 7662977.399999999
$ws.Range("L122").Value = 23971.3329
$ws.Range("M122").Value = -7660527.399999999
$ws.Range("N122").Value = -28871.3329
$ws.Range("H132").Value = 2410.8667
$ws.Range("I132").Value = 2410.8667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7232.6001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4702.6001
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 2120.9143
$ws.Range("I134").Value = 2127.0588
$ws.Range("J134").Value = 1912
$ws.Range("K134").Value = 6381.176399999999
$ws.Range("L134").Value = 5736
$ws.Range("M134").Value = -3846.176399999999
$ws.Range("N134").Value = -10806

# --- Sheet: CUL ---
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H97").Value = 287.5
$ws.Range("I97").Value = 175
$ws.Range("K97").Value = 525
$ws.Range("M97").Value = -29

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 223.81818
$ws.Range("I2").Value = 149.28572
$ws.Range("K2").Value = 149.28572
$ws.Range("M2").Value = -36.28572
$ws.Range("H33").Value = 2000
$ws.Range("J33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("N33").Value = -2504
$ws.Range("H40").Value = 750
$ws.Range("J40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1302
$ws.Range("H44").Value = 25999
$ws.Range("J44").Value = 25999
$ws.Range("L44").Value = 25999
$ws.Range("N44").Value = -27191
$ws.Range("H107").Value = 999.8570999999999
$ws.Range("I107").Value = 998
$ws.Range("J107").Value = 1001.25
$ws.Range("K107").Value = 998
$ws.Range("L107").Value = 1001.25
$ws.Range("M107").Value = 922
$ws.Range("N107").Value = -4841.25
$ws.Range("H109").Value = 49545.453
$ws.Range("J109").Value = 49545.453
$ws.Range("L109").Value = 49545.453
$ws.Range("N109").Value = -51625.453
$ws.Range("H132").Value = 2120
$ws.Range("I132").Value = 1985.7142
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5957.142599999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3427.142599999999
$ws.Range("N132").Value = -17060

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H16").Value = 1779.8572
$ws.Range("I16").Value = 1993.1666
$ws.Range("K16").Value = 1993.1666
$ws.Range("M16").Value = -1823.1666
$ws.Range("H40").Value = 8916.875
$ws.Range("I40").Value = 10709.167
$ws.Range("J40").Value = 3540
$ws.Range("K40").Value = 10709.167
$ws.Range("L40").Value = 3540
$ws.Range("M40").Value = -10573.167
$ws.Range("N40").Value = -3812
$ws.Range("H122").Value = 5151.8
$ws.Range("I122").Value = 2308.182
$ws.Range("K122").Value = 6924.545999999999
$ws.Range("M122").Value = -4474.545999999999
$ws.Range("H136").Value = 3322.5557
$ws.Range("I136").Value = 3322.5557
$ws.Range("K136").Value = 9967.667099999999
$ws.Range("M136").Value = -7417.667099999999

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H81").Value = 12349531
$ws.Range("I81").Value = 2610.8333
$ws.Range("K81").Value = 5221.6666
$ws.Range("M81").Value = -4160.6666
$ws.Range("H84").Value = 12349531
$ws.Range("I84").Value = 2610.8333
$ws.Range("K84").Value = 26108.333
$ws.Range("M84").Value = -20804.333
$ws.Range("H107").Value = 609.375
$ws.Range("I107").Value = 526
$ws.Range("J107").Value = 692.75
$ws.Range("K107").Value = 1578
$ws.Range("L107").Value = 2078.25
$ws.Range("M107").Value = 342
$ws.Range("N107").Value = -5918.25
$ws.Range("H132").Value = 1841.138
$ws.Range("I132").Value = 1841.138
$ws.Range("K132").Value = 5523.414
$ws.Range("M132").Value = -2993.414
$ws.Range("H136").Value = 1355.9
$ws.Range("I136").Value = 761.5454999999999
$ws.Range("J136").Value = 4157.857
$ws.Range("K136").Value = 2284.6365
$ws.Range("L136").Value = 12473.571
$ws.Range("M136").Value = 265.3635000000004
$ws.Range("N136").Value = -17573.571
